# Update NGS analysis stats
# Change the "title" (column I) for the passage-2 HeLa and RPE virus rows
# from "passage 1" to "passage 2" to correctly reflect the sample description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I15").Value = "CVB3 virus from passage 2 of the mutagenized P2 library 1"
$ws.Range("I16").Value = "CVB3 virus from passage 2 of the mutagenized P2 library 2"
$ws.Range("I17").Value = "CVB3 virus from passage 2 of the mutagenized P2 library 3"
$ws.Range("I21").Value = "CVB3 virus from passage 2 of the mutagenized P3 library 1"
$ws.Range("I22").Value = "CVB3 virus from passage 2 of the mutagenized P3 library 2"
$ws.Range("I23").Value = "CVB3 virus from passage 2 of the mutagenized P3 library 3"
$ws.Range("I24").Value = "CVB3 virus from passage 2 of the mutagenized P1 library 1 in RPE cells"
$ws.Range("I25").Value = "CVB3 virus from passage 2 of the mutagenized P1 library 2 in RPE cells"
$ws.Range("I26").Value = "CVB3 virus from passage 2 of the mutagenized P2 library 1 in RPE cells"
$ws.Range("I27").Value = "CVB3 virus from passage 2 of the mutagenized P2 library 2 in RPE cells"
$ws.Range("I28").Value = "CVB3 virus from passage 2 of the mutagenized P2 library 3 in RPE cells"
$ws.Range("I29").Value = "CVB3 virus from passage 2 of the mutagenized P3 library 1 in RPE cells"
$ws.Range("I30").Value = "CVB3 virus from passage 2 of the mutagenized P3 library 2 in RPE cells"
$ws.Range("I31").Value = "CVB3 virus from passage 2 of the mutagenized P3 library 3 in RPE cells"

$ws.Range("I19").Select()
